$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 30 data: group/question label, hours spent, and date
$ws.Range("B30").Value = "Creating and editing Questions"
$ws.Range("C30").Value = 7.25
$ws.Range("D30").Value = 40871

# Copy the date cell formatting/style from D29 onto the new D30 cell
$ws.Range("D29").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect where the user ended up after editing (C31)
$ws.Range("C31").Select()
